# Source diff: cell C10 on the "Rules" worksheet changes its stored
# numeric value from 18 to 1 (serialized upstream as 1.0).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 1
